# Add two newly-collected submissions to the raw data sheet
# ("八位序列号收集收集结果yd5" - the first worksheet), mirroring a fresh
# upload of rows from the collection form. The derived/curated result
# sheet ("八位序列号收集（收集结果）") is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 52: new submission
$ws.Cells.Item(52, 1).Value = "嘻嘻嘻"
$ws.Cells.Item(52, 2).Value = 45923.5852662037
$ws.Cells.Item(52, 2).NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Cells.Item(52, 3).Value = "2a14ed98"
# Column D (QQ number) holds purely-numeric-looking strings that must be
# stored as text, matching every other row in the sheet. Force text
# interpretation via a Text number format, then drop the format again so
# the cell keeps the sheet's default (unstyled) look.
$ws.Cells.Item(52, 4).NumberFormat = "@"
$ws.Cells.Item(52, 4).Value = "3534569125"
$ws.Cells.Item(52, 4).ClearFormats()

# Row 53: new submission
$ws.Cells.Item(53, 1).Value = " "
$ws.Cells.Item(53, 2).Value = 45923.7846296296
$ws.Cells.Item(53, 2).NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Cells.Item(53, 3).Value = "4b24a6b7"
$ws.Cells.Item(53, 4).NumberFormat = "@"
$ws.Cells.Item(53, 4).Value = "2751393486"
$ws.Cells.Item(53, 4).ClearFormats()
